$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) First paragraph: add two trailing spaces to the existing sentence,
#    then append "(This is a change - Version for main branch)" as three
#    separate red-colored runs (matching how Word splits text typed in
#    separate insertion actions).
# ----------------------------------------------------------------------

$p1 = $d.Paragraphs(1)
$full = $p1.Range
$textRange = $d.Range($full.Start, $full.End - 1)
$textRange.Text = "This is a Microsoft word document.  "

$enDash = [char]0x2013

$chunks = @(
    "(This is a change " + $enDash + " Ve",
    "rsion for main branch",
    ")"
)

foreach ($chunk in $chunks) {
    $p1 = $d.Paragraphs(1)
    $endPos = $p1.Range.End
    $insPoint = $d.Range($endPos - 1, $endPos - 1)
    $insPoint.InsertAfter($chunk)

    $p1 = $d.Paragraphs(1)
    $endPos = $p1.Range.End
    $len = $chunk.Length
    $newRun = $d.Range($endPos - 1 - $len, $endPos - 1)
    $newRun.Font.Color = 255
}

# ----------------------------------------------------------------------
# 2) Remove the trailing paragraph that reads
#    "ank God almighty, we are free at last." (the one immediately after
#    "Shall be lifted-nevermore!").
# ----------------------------------------------------------------------

$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$delRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)
$delRange.Delete()
